# Updates Price (D) and Volume(1h) (E) columns for the cryptos list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to stay a text value even when it looks like a number
    # (e.g. "19.54"), matching the sheet's pre-existing text-typed cells,
    # then restore the Normal style so no stray number format lingers.
    $ws.Range($range).NumberFormat = "@"
    $ws.Range($range).Value = $value
    $ws.Range($range).Style = "Normal"
}

$ws.Range("D2").Value = "25.917.95"
$ws.Range("E2").Value = "  -1.36%  "
$ws.Range("D3").Value = "1.637.08"
$ws.Range("E3").Value = "  -0.63%  "
$ws.Range("E4").Value = "  +0.13%  "
Set-TextValue "D5" "215.47"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  -0.99%  "
$ws.Range("E9").Value = "  -0.24%  "
Set-TextValue "D10" "19.54"
$ws.Range("E10").Value = "  -1.85%  "
Set-TextValue "D11" "0.0792"
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("D13").Value = "1.863.81"
$ws.Range("E13").Value = "  -0.63%  "
$ws.Range("D14").Value = "1.609.44"
$ws.Range("E14").Value = "  -2.25%  "
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("D16").Value = "0.0₃0764"
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("D18").Value = "25.938.11"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("E19").Value = "  +0.16%  "
Set-TextValue "D20" "193.09"
$ws.Range("E20").Value = "  -1.37%  "
$ws.Range("E21").Value = "  -2.13%  "
$ws.Range("E22").Value = "  -1.78%  "
$ws.Range("E23").Value = "  -0.86%  "
$ws.Range("E24").Value = "  +4.67%  "
$ws.Range("E25").Value = "  +0.05%  "
Set-TextValue "D26" "143.14"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("E28").Value = "  -0.93%  "
$ws.Range("E29").Value = "  -0.73%  "
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("E31").Value = "  -0.79%  "
$ws.Range("E32").Value = "  -2.20%  "
$ws.Range("E33").Value = "  -0.27%  "
Set-TextValue "D34" "1.54"
$ws.Range("E34").Value = "  -4.28%  "
$ws.Range("E35").Value = "  +1.62%  "
Set-TextValue "D36" "0.901"
$ws.Range("E36").Value = "  -1.44%  "
$ws.Range("D37").Value = "1.133.65"
$ws.Range("E37").Value = "  -0.38%  "
$ws.Range("E38").Value = "  -1.88%  "
$ws.Range("E39").Value = "  -1.56%  "
$ws.Range("E40").Value = "  -0.71%  "
$ws.Range("E41").Value = "  -0.91%  "
Set-TextValue "D42" "99.22"
$ws.Range("E42").Value = "  -1.36%  "
Set-TextValue "D43" "0.797"
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("D44").Value = "1.773.19"
$ws.Range("E44").Value = "  -0.65%  "
$ws.Range("D45").Value = "0.0₆0114"
$ws.Range("E45").Value = "  +2.74%  "
Set-TextValue "D46" "56.56"
$ws.Range("E46").Value = "  -1.15%  "
$ws.Range("E47").Value = "  +2.19%  "
$ws.Range("E48").Value = "  -0.16%  "
Set-TextValue "D49" "7.66"
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("E50").Value = "  -0.89%  "
Set-TextValue "D51" "0.0959"
$ws.Range("E51").Value = "  -1.38%  "
